# AutoCommit_7 января 2024 г. 16:37:12_SibNout2023
# Adds several "5" grade entries to the gradebook and highlights the
# affected cells with a white ("theme 0") fill, mirroring cells that
# were already marked that way elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlThemeColorLight1 (=2) maps to OOXML theme index "0" (the white
# "Background 1" theme color) which is the fill color used for the
# newly-highlighted cells below.
$xlThemeColorLight1 = 2

# --- Row 4 ---------------------------------------------------------
$ws.Range("J4").Value = 5
$ws.Range("J4").Interior.ThemeColor = $xlThemeColorLight1

# --- Row 14 ----------------------------------------------------------
$ws.Range("I14").Value = 5
$ws.Range("I14").Interior.ThemeColor = $xlThemeColorLight1

$ws.Range("J14").Value = 5
$ws.Range("J14").Interior.ThemeColor = $xlThemeColorLight1

# --- Row 21 ----------------------------------------------------------
$ws.Range("E21").Value = 5

$ws.Range("G21").Value = 5
$ws.Range("G21").Interior.ThemeColor = $xlThemeColorLight1

$ws.Range("I21").Value = 5
$ws.Range("I21").Interior.ThemeColor = $xlThemeColorLight1

$ws.Range("M21").Value = 5

# --- Row 25 ----------------------------------------------------------
$ws.Range("J25").Value = 5

# --- Selection ---------------------------------------------------------
$ws.Range("P4").Select() | Out-Null
